$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3   = 6.127000000000001
    14  = 5.504
    16  = 5.401000000000001
    21  = 9.185
    23  = 8.305
    25  = 5.67
    26  = 5.701999999999999
    29  = 5.664
    40  = 9.196
    53  = 6.504
    57  = 5.034000000000001
    59  = 4.704999999999999
    65  = 6.034000000000001
    69  = 5.326
    79  = 5.577000000000001
    83  = 5.702
    91  = 5.601000000000001
    93  = 4.973000000000001
    100 = 5.784
    103 = 5.267
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
